$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.48196983337402
$ws.Range("C3").Value = 17.48800277709961
$ws.Range("C4").Value = 17.56000518798828
$ws.Range("C5").Value = 17.66872406005859
$ws.Range("C6").Value = 16.72506332397461
